$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.096772333333333
$ws.Range("H2").Value = 3.290317
$ws.Range("I2").Value = 0.2426185621302128
$ws.Range("J2").Value = 0.2426185621302128
$ws.Range("M2").Value = 38.45264233333334
$ws.Range("N2").Value = 115.357927
$ws.Range("O2").Value = 0.2975040117664333
$ws.Range("P2").Value = 0.2975040117664332
$ws.Range("Q2").Value = 42.17379425476211
$ws.Range("R2").Value = 379.564148292859
$ws.Range("S2").Value = 0.07217999556274196
$ws.Range("T2").Value = 0.07217999556274195
$ws.Range("G3").Value = 1.096772333333333
$ws.Range("H3").Value = 3.290317
$ws.Range("I3").Value = 0.2426185621302128
$ws.Range("J3").Value = 0.2426185621302128
$ws.Range("O3").Value = 0.3694391181876273
$ws.Range("P3").Value = 0.3694391181876272
$ws.Range("Q3").Value = 52.37122439995166
$ws.Range("R3").Value = 471.341019599565
$ws.Range("S3").Value = 0.08963278764933588
$ws.Range("T3").Value = 0.08963278764933587
$ws.Range("G4").Value = 1.096772333333333
$ws.Range("H4").Value = 3.290317
$ws.Range("I4").Value = 0.2426185621302128
$ws.Range("J4").Value = 0.2426185621302128
$ws.Range("M4").Value = 18.63107466666667
$ws.Range("N4").Value = 55.893224
$ws.Range("O4").Value = 0.1441466469015163
$ws.Range("P4").Value = 0.1441466469015162
$ws.Range("Q4").Value = 20.43404723466755
$ws.Range("R4").Value = 183.906425112008
$ws.Range("S4").Value = 0.03497265220713738
$ws.Range("T4").Value = 0.03497265220713737
$ws.Range("G5").Value = 1.096772333333333
$ws.Range("H5").Value = 3.290317
$ws.Range("I5").Value = 0.2426185621302128
$ws.Range("J5").Value = 0.2426185621302128
$ws.Range("M5").Value = 24.41680433333333
$ws.Range("N5").Value = 73.25041299999999
$ws.Range("O5").Value = 0.1889102231444233
$ws.Range("P5").Value = 0.1889102231444233
$ws.Range("Q5").Value = 26.77967546121344
$ws.Range("R5").Value = 241.017079150921
$ws.Range("S5").Value = 0.04583312671099763
$ws.Range("T5").Value = 0.04583312671099762
$ws.Range("I6").Value = 0.03766810132102297
$ws.Range("J6").Value = 0.03766810132102297
$ws.Range("M6").Value = 38.45264233333334
$ws.Range("N6").Value = 115.357927
$ws.Range("O6").Value = 0.2975040117664333
$ws.Range("P6").Value = 0.2975040117664332
$ws.Range("Q6").Value = 6.547754389162335
$ws.Range("R6").Value = 58.92978950246101
$ws.Range("S6").Value = 0.01120641125862882
$ws.Range("T6").Value = 0.01120641125862882
$ws.Range("I7").Value = 0.03766810132102297
$ws.Range("J7").Value = 0.03766810132102297
$ws.Range("O7").Value = 0.3694391181876273
$ws.Range("P7").Value = 0.3694391181876272
$ws.Range("S7").Value = 0.01391607013584092
$ws.Range("T7").Value = 0.01391607013584092
$ws.Range("I8").Value = 0.03766810132102297
$ws.Range("J8").Value = 0.03766810132102297
$ws.Range("M8").Value = 18.63107466666667
$ws.Range("N8").Value = 55.893224
$ws.Range("O8").Value = 0.1441466469015163
$ws.Range("P8").Value = 0.1441466469015162
$ws.Range("Q8").Value = 3.172518025314667
$ws.Range("R8").Value = 28.55266222783201
$ws.Range("S8").Value = 0.005429730500572036
$ws.Range("T8").Value = 0.005429730500572035
$ws.Range("I9").Value = 0.03766810132102297
$ws.Range("J9").Value = 0.03766810132102297
$ws.Range("M9").Value = 24.41680433333333
$ws.Range("N9").Value = 73.25041299999999
$ws.Range("O9").Value = 0.1889102231444233
$ws.Range("P9").Value = 0.1889102231444233
$ws.Range("Q9").Value = 4.157717858684333
$ws.Range("R9").Value = 37.419460728159
$ws.Range("S9").Value = 0.007115889425981195
$ws.Range("T9").Value = 0.007115889425981193
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5018676666666667
$ws.Range("H10").Value = 1.505603
$ws.Range("I10").Value = 0.1110188577571507
$ws.Range("J10").Value = 0.1110188577571507
$ws.Range("M10").Value = 38.45264233333334
$ws.Range("N10").Value = 115.357927
$ws.Range("O10").Value = 0.2975040117664333
$ws.Range("P10").Value = 0.2975040117664332
$ws.Range("Q10").Value = 19.29813788499789
$ws.Range("R10").Value = 173.683240964981
$ws.Range("S10").Value = 0.03302855556447935
$ws.Range("T10").Value = 0.03302855556447934
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5018676666666667
$ws.Range("H11").Value = 1.505603
$ws.Range("I11").Value = 0.1110188577571507
$ws.Range("J11").Value = 0.1110188577571507
$ws.Range("O11").Value = 0.3694391181876273
$ws.Range("P11").Value = 0.3694391181876272
$ws.Range("Q11").Value = 23.96433917164834
$ws.Range("R11").Value = 215.679052544835
$ws.Range("S11").Value = 0.04101470891199939
$ws.Range("T11").Value = 0.04101470891199937
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5018676666666667
$ws.Range("H12").Value = 1.505603
$ws.Range("I12").Value = 0.1110188577571507
$ws.Range("J12").Value = 0.1110188577571507
$ws.Range("M12").Value = 18.63107466666667
$ws.Range("N12").Value = 55.893224
$ws.Range("O12").Value = 0.1441466469015163
$ws.Range("P12").Value = 0.1441466469015162
$ws.Range("Q12").Value = 9.350333970452445
$ws.Range("R12").Value = 84.15300573407201
$ws.Range("S12").Value = 0.01600299608852967
$ws.Range("T12").Value = 0.01600299608852966
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5018676666666667
$ws.Range("H13").Value = 1.505603
$ws.Range("I13").Value = 0.1110188577571507
$ws.Range("J13").Value = 0.1110188577571507
$ws.Range("M13").Value = 24.41680433333333
$ws.Range("N13").Value = 73.25041299999999
$ws.Range("O13").Value = 0.1889102231444233
$ws.Range("P13").Value = 0.1889102231444233
$ws.Range("Q13").Value = 12.25400461822656
$ws.Range("R13").Value = 110.286041564039
$ws.Range("S13").Value = 0.02097259719214233
$ws.Range("T13").Value = 0.02097259719214233
$ws.Range("G14").Value = 2.751641333333334
$ws.Range("H14").Value = 8.254924000000001
$ws.Range("I14").Value = 0.6086944787916135
$ws.Range("J14").Value = 0.6086944787916135
$ws.Range("M14").Value = 38.45264233333334
$ws.Range("N14").Value = 115.357927
$ws.Range("O14").Value = 0.2975040117664333
$ws.Range("P14").Value = 0.2975040117664332
$ws.Range("Q14").Value = 105.8078800202831
$ws.Range("R14").Value = 952.2709201825481
$ws.Range("S14").Value = 0.1810890493805832
$ws.Range("T14").Value = 0.1810890493805831
$ws.Range("G15").Value = 2.751641333333334
$ws.Range("H15").Value = 8.254924000000001
$ws.Range("I15").Value = 0.6086944787916135
$ws.Range("J15").Value = 0.6086944787916135
$ws.Range("O15").Value = 0.3694391181876273
$ws.Range("P15").Value = 0.3694391181876272
$ws.Range("Q15").Value = 131.3917404336867
$ws.Range("R15").Value = 1182.52566390318
$ws.Range("S15").Value = 0.2248755514904511
$ws.Range("T15").Value = 0.224875551490451
$ws.Range("G16").Value = 2.751641333333334
$ws.Range("H16").Value = 8.254924000000001
$ws.Range("I16").Value = 0.6086944787916135
$ws.Range("J16").Value = 0.6086944787916135
$ws.Range("M16").Value = 18.63107466666667
$ws.Range("N16").Value = 55.893224
$ws.Range("O16").Value = 0.1441466469015163
$ws.Range("P16").Value = 0.1441466469015162
$ws.Range("Q16").Value = 51.26603513721956
$ws.Range("R16").Value = 461.3943162349761
$ws.Range("S16").Value = 0.0877412681052772
$ws.Range("T16").Value = 0.08774126810527719
$ws.Range("G17").Value = 2.751641333333334
$ws.Range("H17").Value = 8.254924000000001
$ws.Range("I17").Value = 0.6086944787916135
$ws.Range("J17").Value = 0.6086944787916135
$ws.Range("M17").Value = 24.41680433333333
$ws.Range("N17").Value = 73.25041299999999
$ws.Range("O17").Value = 0.1889102231444233
$ws.Range("P17").Value = 0.1889102231444233
$ws.Range("Q17").Value = 67.18628803151245
$ws.Range("R17").Value = 604.676592283612
$ws.Range("S17").Value = 0.1149886098153022
$ws.Range("T17").Value = 0.1149886098153021
